$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: (empty) -> false
# Typing "false"/"true" directly gets auto-converted to a Boolean by Excel,
# so force literal text via a leading apostrophe, then restore the
# original (non quote-prefixed) cell formatting by copying it from a
# neighboring cell that already carries the correct style.
$ws.Range("B7").Value = "'false"
$ws.Range("B13").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date: updated timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Case Sensitive: (empty) -> true
$ws.Range("B15").Value = "'true"
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$excel.CutCopyMode = $false
